$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "MuSCs" target-cluster label (used in column D) to "Resolving-Mac"
# and refresh the TPM-derived statistics for rows 2-10 to the updated values.

$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("D10").Value = "Resolving-Mac"

# Row 2
$ws.Range("G2").Value = 0.06493566666666667
$ws.Range("H2").Value = 0.194807
$ws.Range("I2").Value = 0.05459895593058446
$ws.Range("J2").Value = 0.05459895593058446
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.058783666666667
$ws.Range("N2").Value = 3.176351
$ws.Range("O2").Value = 0.5014862149947701
$ws.Range("P2").Value = 0.5014862149947702
$ws.Range("Q2").Value = 0.06875282325077779
$ws.Range("R2").Value = 0.6187754092570001
$ws.Range("S2").Value = 0.02738062375229506
$ws.Range("T2").Value = 0.02738062375229506

# Row 3
$ws.Range("G3").Value = 0.06493566666666667
$ws.Range("H3").Value = 0.194807
$ws.Range("I3").Value = 0.05459895593058446
$ws.Range("J3").Value = 0.05459895593058446
$ws.Range("O3").Value = 0.475327031240749
$ws.Range("P3").Value = 0.4753270312407492
$ws.Range("Q3").Value = 0.06516644802600001
$ws.Range("R3").Value = 0.586498032234
$ws.Range("S3").Value = 0.0259523596313292
$ws.Range("T3").Value = 0.02595235963132921

# Row 4
$ws.Range("G4").Value = 0.06493566666666667
$ws.Range("H4").Value = 0.194807
$ws.Range("I4").Value = 0.05459895593058446
$ws.Range("J4").Value = 0.05459895593058446
$ws.Range("M4").Value = 0.048954
$ws.Range("N4").Value = 0.146862
$ws.Range("O4").Value = 0.02318675376448066
$ws.Range("P4").Value = 0.02318675376448067
$ws.Range("Q4").Value = 0.003178860626
$ws.Range("R4").Value = 0.028609745634
$ws.Range("S4").Value = 0.001265972546960193
$ws.Range("T4").Value = 0.001265972546960193

# Row 5
$ws.Range("I5").Value = 0.8534208043695533
$ws.Range("J5").Value = 0.8534208043695533
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.058783666666667
$ws.Range("N5").Value = 3.176351
$ws.Range("O5").Value = 0.5014862149947701
$ws.Range("P5").Value = 0.5014862149947702
$ws.Range("Q5").Value = 1.074655892613667
$ws.Range("R5").Value = 9.671903033523002
$ws.Range("S5").Value = 0.4279787689810794
$ws.Range("T5").Value = 0.4279787689810795

# Row 6
$ws.Range("I6").Value = 0.8534208043695533
$ws.Range("J6").Value = 0.8534208043695533
$ws.Range("O6").Value = 0.475327031240749
$ws.Range("P6").Value = 0.4753270312407492
$ws.Range("S6").Value = 0.4056539773400718
$ws.Range("T6").Value = 0.4056539773400719

# Row 7
$ws.Range("I7").Value = 0.8534208043695533
$ws.Range("J7").Value = 0.8534208043695533
$ws.Range("M7").Value = 0.048954
$ws.Range("N7").Value = 0.146862
$ws.Range("O7").Value = 0.02318675376448066
$ws.Range("P7").Value = 0.02318675376448067
$ws.Range("Q7").Value = 0.049687869414
$ws.Range("R7").Value = 0.447190824726
$ws.Range("S7").Value = 0.01978805804840186
$ws.Range("T7").Value = 0.01978805804840186

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.109394
$ws.Range("H8").Value = 0.328182
$ws.Range("I8").Value = 0.09198023969986227
$ws.Range("J8").Value = 0.09198023969986227
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.058783666666667
$ws.Range("N8").Value = 3.176351
$ws.Range("O8").Value = 0.5014862149947701
$ws.Range("P8").Value = 0.5014862149947702
$ws.Range("Q8").Value = 0.1158245804313333
$ws.Range("R8").Value = 1.042421223882
$ws.Range("S8").Value = 0.04612682226139562
$ws.Range("T8").Value = 0.04612682226139563

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.109394
$ws.Range("H9").Value = 0.328182
$ws.Range("I9").Value = 0.09198023969986227
$ws.Range("J9").Value = 0.09198023969986227
$ws.Range("O9").Value = 0.475327031240749
$ws.Range("P9").Value = 0.4753270312407492
$ws.Range("Q9").Value = 0.109782786276
$ws.Range("R9").Value = 0.9880450764839999
$ws.Range("S9").Value = 0.04372069426934802
$ws.Range("T9").Value = 0.04372069426934803

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.109394
$ws.Range("H10").Value = 0.328182
$ws.Range("I10").Value = 0.09198023969986227
$ws.Range("J10").Value = 0.09198023969986227
$ws.Range("M10").Value = 0.048954
$ws.Range("N10").Value = 0.146862
$ws.Range("O10").Value = 0.02318675376448066
$ws.Range("P10").Value = 0.02318675376448067
$ws.Range("Q10").Value = 0.005355273875999999
$ws.Range("R10").Value = 0.04819746488399999
$ws.Range("S10").Value = 0.002132723169118615
$ws.Range("T10").Value = 0.002132723169118616
